$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Decision country samples with few elections (Comoros) and early elections (Kazakhstan)
$ws.Range("A62").Value = "Comoros"
$ws.Range("B62").Value = 1
$ws.Range("C62").Value = 0
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0

$ws.Range("A63").Value = "Kazakhstan"
$ws.Range("B63").Value = 0
$ws.Range("C63").Value = 0
$ws.Range("D63").Value = 1
$ws.Range("E63").Value = 0

# Move the selection/view to reflect the new end of the data range
$ws.Range("D67").Select()
